# Update countries & provincias Spain
#
# This reshuffles a handful of country rows (their labels moved around in
# the shared-strings table while the stats kept up with a newer data pull)
# and refreshes several countries' case/death counters.
#
# Net effect per touched row (country name + full B:H stats), derived from
# the canonical-XML diff:
#   18  Suiza               -> stats refreshed (E,G,H)
#   44  Malasia             -> stats refreshed (B..H)
#   65  Barein -> Kuwait     (Kuwait gets refreshed stats)
#   66  Kuwait -> Barein     (keeps Barein's prior stats)
#   70  Uzbekistan          -> stats refreshed (D,E)
#  123  El Salvador -> Mali  (Mali gets refreshed stats)
#  124  Islas Feroe -> El Salvador (keeps prior stats)
#  125  Mali -> Islas Feroe  (keeps prior stats)
#  160  Macao -> Guinea-Bisau (Guinea-Bisau gets refreshed stats)
#  161  Guinea-Bisau -> Macao (keeps prior stats)
#  169  Guam -> Republica del Chad (gets refreshed stats)
#  170  Maldivas -> Guam          (keeps prior stats)
#  171  Mongolia -> Maldivas      (keeps prior stats)
#  172  Nepal -> Mongolia         (keeps prior stats)
#  173  Republica del Chad -> Nepal (keeps prior stats)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($r, $country, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($r, 1).Value = $country
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
}

# Row 18 - Suiza: refreshed stats (name unchanged)
Set-Row 18 "Suiza" 27078 0 16400 9342 386 9 1336

# Row 44 - Malasia: refreshed stats (name unchanged)
Set-Row 44 "Malasia" 5305 54 3102 2115 49 2 88

# Rows 65/66 - Kuwait moves ahead of Barein; Kuwait's stats refreshed
Set-Row 65 "Kuwait" 1751 93 280 1465 34 1 6
Set-Row 66 "Barein" 1744 4 726 1011 3 0 7

# Row 70 - Uzbekistan: refreshed stats (name unchanged)
Set-Row 70 "Uzbekistan" 1450 45 161 1285 8 0 4

# Rows 123/124/125 - Mali moves ahead of El Salvador / Islas Feroe; Mali's stats refreshed
Set-Row 123 "Mali" 190 19 34 143 0 0 13
Set-Row 124 "El Salvador" 190 13 43 140 2 0 7
Set-Row 125 "Islas Feroe" 184 0 173 11 0 0 0

# Rows 160/161 - Guinea-Bisau moves ahead of Macao; Guinea-Bisau's stats refreshed
Set-Row 160 "Guinea-Bisau" 46 3 0 46 0 0 0
Set-Row 161 "Macao" 45 0 16 29 1 0 0

# Rows 169-173 - Republica del Chad moves ahead of Guam/Maldivas/Mongolia/Nepal;
# Republica del Chad's stats refreshed, the other four keep their prior stats
# shifted down by one row.
Set-Row 169 "Republica del Chad" 33 6 8 25 0 0 0
Set-Row 170 "Guam" 32 0 0 31 0 0 1
Set-Row 171 "Maldivas" 32 3 16 16 0 0 0
Set-Row 172 "Mongolia" 31 0 5 26 0 0 0
Set-Row 173 "Nepal" 30 0 2 28 0 0 0
